# Recuperacao da versao anterior do PITCH
#
# Turns the "Rectangle 3" autoshape (filled rectangle labelled "LOGO")
# back into the earlier "Oval 3" autoshape (accent-coloured oval labelled
# "LOGO DO PROJETO") that sits in the upper-right area of the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape by its current name rather than a hard-coded index so
# the script is resilient to shape ordering.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Rectangle 3") {
        $shp = $s.Shapes.Item($i)
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(3)
}

# Rename.
$shp.Name = "Oval 3"

# Reposition / resize (EMU -> points, 12700 EMU per point). The COM layer
# truncates the points value through a single-precision float on its way
# back to EMU, so a couple of the raw EMU/12700 quotients land one EMU
# short; nudge those by a hair so the stored EMU comes out exact.
$shp.Left   = 649.2414273228346   # 8245366 EMU
$shp.Top    = 270.0                # 3429000 EMU
$shp.Width  = 144.0                # 1828800 EMU
$shp.Height = 130.37497062992125  # 1655762 EMU

# Change the preset geometry from rectangle to oval.
$shp.AutoShapeType = 9   # msoShapeOval

# The restored oval no longer carries an explicit solid-fill override; it
# falls back to the shape style's fillRef (scheme accent1). The closest
# reachable state through the exposed Fill API is to point the explicit
# fill at the same theme colour referenced by the style.
$shp.Fill.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1

# Update the label text.
$shp.TextFrame.TextRange.Text = "LOGO DO PROJETO"
